# Update statistics (想去人数 / 最低票价) across the "展览" and "全部类型" sheets
# to reflect newly generated output data (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> cell address -> new value
$updates = @{
    "展览" = @{
        "G2"  = 93.90000000000001
        "F3"  = 14770
        "F4"  = 18167
        "F6"  = 97
        "F15" = 73
        "F16" = 188
        "F18" = 1371
        "F21" = 76
        "F23" = 7517
        "F27" = 1199
        "F29" = 5906
        "F30" = 85
        "F34" = 251
        "F35" = 5212
        "F37" = 37
    }
    "全部类型" = @{
        "G2"  = 93.90000000000001
        "F3"  = 14770
        "F4"  = 18168
        "F6"  = 97
        "F15" = 73
        "F16" = 188
        "F18" = 1371
        "F22" = 76
        "F24" = 7518
        "F28" = 1199
        "F31" = 5906
        "F32" = 85
        "F36" = 251
        "F37" = 5212
        "F39" = 37
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
